$d = $word.ActiveDocument

function Split-RunAt($rangeObj) {
    # Toggling a character formatting property on/off at the exact same
    # range forces Word to break runs at that range's boundaries without
    # altering the visible formatting.
    $rangeObj.Bold = 1
    $rangeObj.Bold = 0
}

# ---------------------------------------------------------------
# Edit 1: "NIR spectra of the milled barley samples." (Xcal paragraph)
#   -> "NIR spectra of the " + "tablets" + "."
# ---------------------------------------------------------------
$full = $d.Content.Text
$needle = "NIR spectra of the milled barley samples."
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $prefix = "NIR spectra of the "
    $word1  = "milled barley samples"
    $midStart = $idx + $prefix.Length
    $midEnd   = $midStart + $word1.Length

    $rMid = $d.Range($midStart, $midEnd)
    $rMid.Text = "tablets"

    $rSplit = $d.Range($midStart, $midStart + "tablets".Length)
    Split-RunAt $rSplit
}

# ---------------------------------------------------------------
# Edit 2: "X" + "v" + "al" (three bold runs) -> single bold run "Xval"
# ---------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("Xval")
if ($idx -ge 0) {
    $find = $d.Content.Find
    $find.Text = "Xval"
    $find.Replacement.Text = "Xval"
    $find.Execute("Xval", $false, $false, $false, $false, $false, $true, 1, $false, "Xval", 2) | Out-Null
}

# ---------------------------------------------------------------
# Edit 3: "NIR spectra of the milled barley samples (test set samples)."
#   -> "NIR spectra of the " + "tablets" + " " + "(test set samples)" + "."
# ---------------------------------------------------------------
$full = $d.Content.Text
$needle = "NIR spectra of the milled barley samples ("
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $prefix = "NIR spectra of the "
    $word1  = "milled barley samples"
    $midStart = $idx + $prefix.Length
    $midEnd   = $midStart + $word1.Length

    $rMid = $d.Range($midStart, $midEnd)
    $rMid.Text = "tablets"

    $rSplit1 = $d.Range($midStart, $midStart + "tablets".Length)
    Split-RunAt $rSplit1

    $afterStart = $midStart + "tablets".Length
    $spaceEnd = $afterStart + 1
    $rSplit2 = $d.Range($afterStart, $spaceEnd)
    Split-RunAt $rSplit2

    $parenStart = $spaceEnd
    $parenLen = "(test set samples)".Length
    $parenEnd = $parenStart + $parenLen
    $rSplit3 = $d.Range($parenStart, $parenEnd)
    Split-RunAt $rSplit3
}

# ---------------------------------------------------------------
# Edit 4: remove the leftover "_GoBack" bookmark
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
